# #5: cash & deposit done
# Fill out the "存款" (deposit) sheet with bank / deposit_type / currency
# headers and the shared metadata columns (property_category, category,
# date, legislator_name, legislator_id, source_file, index) that the other
# property sheets already carry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# ---- Row 1: proper column headers (was a copy of row 2's data) ----
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# ---- Rows 2-10: add the new metadata columns G..M ----
# (columns A-F already hold the correct bank/deposit_type/currency/owner/total
# data and are left untouched)
$indices = @(50, 51, 52, 53, 54, 55, 56, 57, 58)

for ($i = 0; $i -lt $indices.Length; $i++) {
    $r = $i + 2
    $idx = $indices[$i]

    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    $ws.Cells.Item($r, 9).Value = "2011-11-21"
    $ws.Cells.Item($r, 10).Value = "孫大千"
    $ws.Cells.Item($r, 11).Value = 919
    $ws.Cells.Item($r, 12).Value = "tmpc6841"
    $ws.Cells.Item($r, 13).Value = $idx
}
